# session-25.pptx - "Todays Assignment" slide (slide 5):
#   - merge the two runs of the "Complete Ch.10.1 through 10.2 / on Ajax"
#     paragraph into a single run of text
#   - add a new paragraph "Review associated AJAX slides" right after it

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# The content placeholder currently holds 3 paragraphs (1-based Characters
# indexing):
#   1) "Assignment (before next class):"              (31 chars, starts at 1)
#   2) "Project 4: Contact Manager (client)"           (35 chars, starts at 33)
#   3) "Complete Ch.10.1 through 10.2 " + "on Ajax"     (37 chars, starts at 69)
# Paragraphs are separated by a single carriage-return character, so
# paragraph 3's text starts right after "1 + 31 + 1 + 35 + 1" characters.
$para3Start = 1 + 31 + 1 + 35 + 1
$para3Len = 37

# Re-type paragraph 3 as one run (this merges the two existing runs into one).
$para3 = $tr.Characters($para3Start, $para3Len)
$para3.Text = "Complete Ch.10.1 through 10.2 on Ajax"

# Append the new, fourth paragraph after the (now merged) paragraph three.
$tr = $shape.TextFrame.TextRange
$null = $tr.InsertAfter([char]13 + "Review associated AJAX slides")
